$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "SEC_Comm": add the new "WIND_OFF" commodity row (row 8, previously
# a blank templated row right under the existing SOLAR commodity row).
# ---------------------------------------------------------------------------
$secComm = $wb.Worksheets.Item("SEC_Comm")
$secComm.Range("B8").Value = "NRG"
$secComm.Range("C8").Value = "WIND_OFF"
$secComm.Range("D8").Value = "Wind Offshore Energy"
$secComm.Range("E8").Value = "PJ"
$secComm.Range("G8").Value = "SEASON"

# ---------------------------------------------------------------------------
# Sheet "SEC_Processes": add the new "MIN_OFFSHORE" mining process (row 9)
# and the new "NEW_OFFSHORE_PP" power-plant process (row 10) - both were
# blank templated rows below the existing SOLAR / NEW_PV_PP rows.
# ---------------------------------------------------------------------------
$secProcesses = $wb.Worksheets.Item("SEC_Processes")
$secProcesses.Range("B9").Value = "MIN"
$secProcesses.Range("C9").Value = "PL"
$secProcesses.Range("D9").Value = "MIN_OFFSHORE"
$secProcesses.Range("E9").Value = "Offshore Wind Energy Supply"
$secProcesses.Range("F9").Value = "PJ"
$secProcesses.Range("G9").Value = "PJ/a"
$secProcesses.Range("H9").Value = "SEASON"

$secProcesses.Range("B10").Value = "ELE"
$secProcesses.Range("C10").Value = "PL"
$secProcesses.Range("D10").Value = "NEW_OFFSHORE_PP"
$secProcesses.Range("E10").Value = "New Offshore Power Plant"
$secProcesses.Range("F10").Value = "PJ"
$secProcesses.Range("G10").Value = "GW"
$secProcesses.Range("H10").Value = "DAYNITE"

# ---------------------------------------------------------------------------
# Sheet "MIN_IMP": add the mining-cost row (row 9) for MIN_OFFSHORE / WIND_OFF
# ---------------------------------------------------------------------------
$minImp = $wb.Worksheets.Item("MIN_IMP")
$minImp.Range("B9").Formula = "=SEC_Processes!D9"
$minImp.Range("D9").Formula = "=SEC_Comm!C8"
$minImp.Range("E9").Value = 0.001
$minImp.Range("F9").Value = 0.002
$minImp.Range("G9").Value = 2025
$minImp.Range("H9").Value = 100

# ---------------------------------------------------------------------------
# Sheet "PP": add the new offshore power plant row (row 9)
# ---------------------------------------------------------------------------
$pp = $wb.Worksheets.Item("PP")
$pp.Range("B9").Formula = "=SEC_Processes!D10"
$pp.Range("C9").Value = "New Offshore Power Plant"
$pp.Range("D9").Formula = "=SEC_Comm!C8"
$pp.Range("E9").Formula = "=SEC_Comm!C27"
$pp.Range("F9").Value = 1
$pp.Range("G9").Value = 31.536
$pp.Range("H9").Value = 0.4
$pp.Range("I9").Value = 1
$pp.Range("K9").Value = 3000
$pp.Range("L9").Value = 2025
$pp.Range("M9").Value = 25
